$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 599.2857  # H6
$ws.Cells.Item(6, 9).Value = 65  # I6
$ws.Cells.Item(6, 11).Value = 195  # K6
$ws.Cells.Item(6, 13).Value = -83  # M6
$ws.Cells.Item(17, 8).Value = 1679.625  # H17
$ws.Cells.Item(17, 10).Value = 1634.2  # J17
$ws.Cells.Item(17, 12).Value = 4902.6  # L17
$ws.Cells.Item(17, 14).Value = -5238.6  # N17
$ws.Cells.Item(62, 8).Value = 2909.7144  # H62
$ws.Cells.Item(62, 9).Value = 2342  # I62
$ws.Cells.Item(62, 11).Value = 2342  # K62
$ws.Cells.Item(62, 13).Value = -1718  # M62
$ws.Cells.Item(65, 8).Value = 2909.7144  # H65
$ws.Cells.Item(65, 9).Value = 2342  # I65
$ws.Cells.Item(65, 11).Value = 11710  # K65
$ws.Cells.Item(65, 13).Value = -8590  # M65
$ws.Cells.Item(141, 8).Value = 4981.96  # H141
$ws.Cells.Item(141, 9).Value = 3890.4285  # I141
$ws.Cells.Item(141, 11).Value = 11671.2855  # K141
$ws.Cells.Item(141, 13).Value = -6491.2855  # M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1842.5  # H61
$ws.Cells.Item(61, 9).Value = 1842.5  # I61
$ws.Cells.Item(61, 11).Value = 1842.5  # K61
$ws.Cells.Item(61, 13).Value = -1630.5  # M61
$ws.Cells.Item(122, 8).Value = 5357.2856  # H122
$ws.Cells.Item(122, 10).Value = 5168.5835  # J122
$ws.Cells.Item(122, 12).Value = 15505.7505  # L122
$ws.Cells.Item(122, 14).Value = -20405.7505  # N122
$ws.Cells.Item(132, 8).Value = 1564.5  # H132
$ws.Cells.Item(132, 9).Value = 1592.1428  # I132
$ws.Cells.Item(132, 11).Value = 4776.428400000001  # K132
$ws.Cells.Item(132, 13).Value = -2246.428400000001  # M132
$ws.Cells.Item(136, 8).Value = 1842.5  # H136
$ws.Cells.Item(136, 9).Value = 1842.5  # I136
$ws.Cells.Item(136, 11).Value = 5527.5  # K136
$ws.Cells.Item(136, 13).Value = -2977.5  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3449.5  # H99
$ws.Cells.Item(99, 9).Value = 3819.8  # I99
$ws.Cells.Item(99, 10).Value = 1598  # J99
$ws.Cells.Item(99, 11).Value = 3819.8  # K99
$ws.Cells.Item(99, 12).Value = 1598  # L99
$ws.Cells.Item(99, 13).Value = -2321.8  # M99
$ws.Cells.Item(99, 14).Value = -4594  # N99
$ws.Cells.Item(107, 8).Value = 5607.409  # H107
$ws.Cells.Item(107, 9).Value = 5108.5  # I107
$ws.Cells.Item(107, 10).Value = 7852.5  # J107
$ws.Cells.Item(107, 11).Value = 5108.5  # K107
$ws.Cells.Item(107, 12).Value = 7852.5  # L107
$ws.Cells.Item(107, 13).Value = -3188.5  # M107
$ws.Cells.Item(107, 14).Value = -11692.5  # N107
$ws.Cells.Item(134, 8).Value = 2799  # H134
$ws.Cells.Item(134, 9).Value = 2035.5294  # I134
$ws.Cells.Item(134, 10).Value = 5394.8  # J134
$ws.Cells.Item(134, 11).Value = 6106.5882  # K134
$ws.Cells.Item(134, 12).Value = 16184.4  # L134
$ws.Cells.Item(134, 13).Value = -3571.5882  # M134
$ws.Cells.Item(134, 14).Value = -21254.4  # N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2692.9048  # H16
$ws.Cells.Item(16, 9).Value = 2473.7646  # I16
$ws.Cells.Item(16, 11).Value = 2473.7646  # K16
$ws.Cells.Item(16, 13).Value = -2186.7646  # M16
$ws.Cells.Item(31, 8).Value = 1278.6052  # H31
$ws.Cells.Item(31, 9).Value = 1176.5483  # I31
$ws.Cells.Item(31, 10).Value = 1730.5714  # J31
$ws.Cells.Item(31, 11).Value = 1176.5483  # K31
$ws.Cells.Item(31, 12).Value = 1730.5714  # L31
$ws.Cells.Item(31, 13).Value = -881.5482999999999  # M31
$ws.Cells.Item(31, 14).Value = -2320.5714  # N31
$ws.Cells.Item(34, 8).Value = 1278.6052  # H34
$ws.Cells.Item(34, 9).Value = 1176.5483  # I34
$ws.Cells.Item(34, 10).Value = 1730.5714  # J34
$ws.Cells.Item(34, 11).Value = 1176.5483  # K34
$ws.Cells.Item(34, 12).Value = 1730.5714  # L34
$ws.Cells.Item(34, 13).Value = -974.5482999999999  # M34
$ws.Cells.Item(34, 14).Value = -2134.5714  # N34
$ws.Cells.Item(36, 8).Value = 0  # H36
$ws.Cells.Item(36, 9).Value = 0  # I36
$ws.Cells.Item(36, 11).Value = 0  # K36
$ws.Cells.Item(36, 13).ClearContents()  # M36
$ws.Cells.Item(40, 8).Value = 0  # H40
$ws.Cells.Item(40, 9).Value = 0  # I40
$ws.Cells.Item(40, 11).Value = 0  # K40
$ws.Cells.Item(40, 13).ClearContents()  # M40
$ws.Cells.Item(94, 8).Value = 12768.889  # H94
$ws.Cells.Item(94, 10).Value = 1985.1428  # J94
$ws.Cells.Item(94, 12).Value = 1985.1428  # L94
$ws.Cells.Item(94, 14).Value = -2887.1428  # N94
$ws.Cells.Item(105, 8).Value = 658.3333  # H105
$ws.Cells.Item(105, 9).Value = 658.3333  # I105
$ws.Cells.Item(105, 11).Value = 658.3333  # K105
$ws.Cells.Item(105, 13).Value = 1088.6667  # M105
$ws.Cells.Item(113, 8).Value = 2692.9048  # H113
$ws.Cells.Item(113, 9).Value = 2473.7646  # I113
$ws.Cells.Item(113, 11).Value = 2473.7646  # K113
$ws.Cells.Item(113, 13).Value = -303.7646  # M113
$ws.Cells.Item(132, 8).Value = 2698.1614  # H132
$ws.Cells.Item(132, 9).Value = 2008.0476  # I132
$ws.Cells.Item(132, 10).Value = 4147.4  # J132
$ws.Cells.Item(132, 11).Value = 6024.142800000001  # K132
$ws.Cells.Item(132, 12).Value = 12442.2  # L132
$ws.Cells.Item(132, 13).Value = -3494.142800000001  # M132
$ws.Cells.Item(132, 14).Value = -17502.2  # N132
$ws.Cells.Item(134, 8).Value = 3598.2144  # H134
$ws.Cells.Item(134, 9).Value = 3489.4211  # I134
$ws.Cells.Item(134, 11).Value = 10468.2633  # K134
$ws.Cells.Item(134, 13).Value = -7933.263300000001  # M134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 69.25  # H2
$ws.Cells.Item(2, 9).Value = 69.5  # I2
$ws.Cells.Item(2, 11).Value = 417  # K2
$ws.Cells.Item(2, 13).Value = -304  # M2
$ws.Cells.Item(13, 8).Value = 199.05556  # H13
$ws.Cells.Item(13, 9).Value = 247.3077  # I13
$ws.Cells.Item(13, 10).Value = 73.59999999999999  # J13
$ws.Cells.Item(13, 11).Value = 741.9231  # K13
$ws.Cells.Item(13, 12).Value = 220.8  # L13
$ws.Cells.Item(13, 13).Value = -573.9231  # M13
$ws.Cells.Item(13, 14).Value = -556.8  # N13
$ws.Cells.Item(16, 8).Value = 582.7143  # H16
$ws.Cells.Item(16, 9).Value = 263.16666  # I16
$ws.Cells.Item(16, 10).Value = 2500  # J16
$ws.Cells.Item(16, 11).Value = 789.4999799999999  # K16
$ws.Cells.Item(16, 12).Value = 7500  # L16
$ws.Cells.Item(16, 13).Value = -616.4999799999999  # M16
$ws.Cells.Item(16, 14).Value = -7846  # N16
$ws.Cells.Item(56, 8).Value = 7289.1665  # H56
$ws.Cells.Item(56, 9).Value = 7289.1665  # I56
$ws.Cells.Item(56, 11).Value = 7289.1665  # K56
$ws.Cells.Item(56, 13).Value = -6759.1665  # M56
$ws.Cells.Item(140, 8).Value = 22733184  # H140
$ws.Cells.Item(140, 9).Value = 50003410  # I140
$ws.Cells.Item(140, 11).Value = 150010230  # K140
$ws.Cells.Item(140, 13).Value = -150005050  # M140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(38, 8).Value = 17578.8  # H38
$ws.Cells.Item(38, 10).Value = 17578.8  # J38
$ws.Cells.Item(38, 12).Value = 17578.8  # L38
$ws.Cells.Item(38, 14).Value = -18504.8  # N38
$ws.Cells.Item(40, 8).Value = 25000  # H40
$ws.Cells.Item(40, 10).Value = 25000  # J40
$ws.Cells.Item(40, 12).Value = 25000  # L40
$ws.Cells.Item(40, 14).Value = -25302  # N40
$ws.Cells.Item(47, 8).Value = 0  # H47
$ws.Cells.Item(47, 10).Value = 0  # J47
$ws.Cells.Item(47, 12).Value = 0  # L47
$ws.Cells.Item(47, 14).ClearContents()  # N47
$ws.Cells.Item(132, 8).Value = 8533.6  # H132
$ws.Cells.Item(132, 9).Value = 8293.212  # I132
$ws.Cells.Item(132, 11).Value = 24879.636  # K132
$ws.Cells.Item(132, 13).Value = -22349.636  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1991.6774  # H132
$ws.Cells.Item(132, 9).Value = 1724.7333  # I132
$ws.Cells.Item(132, 10).Value = 10000  # J132
$ws.Cells.Item(132, 11).Value = 5174.199900000001  # K132
$ws.Cells.Item(132, 12).Value = 30000  # L132
$ws.Cells.Item(132, 13).Value = -2644.199900000001  # M132
$ws.Cells.Item(132, 14).Value = -35060  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 766.7895  # H107
$ws.Cells.Item(107, 9).Value = 849.8182  # I107
$ws.Cells.Item(107, 11).Value = 2549.4546  # K107
$ws.Cells.Item(107, 13).Value = -629.4546  # M107
